$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.326.41'
$ws.Range('D3').Value = '1.628.64'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('D5').Value = "'1.002"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'297.47"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.97%  '
$ws.Range('D7').Value = "'0.3761"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.27%  '
$ws.Range('D8').Value = "'50.04"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.34%  '
$ws.Range('D9').Value = "'0.3505"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.18%  '
$ws.Range('D10').Value = "'0.08013"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('D11').Value = "'1.199"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.67%  '
$ws.Range('D12').Value = "'1.003"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').Value = "'21.79"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.74%  '
$ws.Range('D14').Value = "'6.267"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.80%  '
$ws.Range('D15').Value = "'7.204"
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Value = "'0.00001188"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.71%  '
$ws.Range('D17').Value = '1.630.72'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('D18').Value = "'95.03"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.55%  '
$ws.Range('D19').Value = "'0.06926"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').Value = "'6.628"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('D21').Value = "'17.24"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.64%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').Value = "'12.21"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.08%  '
$ws.Range('D24').Value = '23.352.68'
$ws.Range('E24').Value = '  -1.42%  '
$ws.Range('D25').Value = "'2.436"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('D26').Value = "'2.909"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.81%  '
$ws.Range('D27').Value = "'20.68"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.76%  '
$ws.Range('D28').Value = "'151.60"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.35%  '
$ws.Range('D29').Value = "'5.172"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').Value = "'131.47"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.65%  '
$ws.Range('D31').Value = '1.813.51'
$ws.Range('E31').Value = '  -1.53%  '
$ws.Range('D32').Value = "'6.760"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.01%  '
$ws.Range('E33').Value = '  -4.42%  '
$ws.Range('D34').Value = "'11.13"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.01%  '
$ws.Range('D35').Value = "'0.9617"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -9.82%  '
$ws.Range('D36').Value = "'0.02676"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.03%  '
$ws.Range('D37').Value = "'0.08689"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('D38').Value = "'0.2410"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.88%  '
$ws.Range('D39').Value = "'5.814"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.83%  '
$ws.Range('D40').Value = "'0.06745"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.08%  '
$ws.Range('E41').Value = '  -2.26%  '
$ws.Range('D42').Value = "'0.6781"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('D43').Value = "'1.289"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.25%  '
$ws.Range('D44').Value = "'15.37"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.88%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').Value = "'0.6280"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.71%  '
$ws.Range('D47').Value = "'2.223"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.07%  '
$ws.Range('D48').Value = "'3.884"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.31%  '
$ws.Range('D49').Value = "'0.07656"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.79%  '
$ws.Range('D50').Value = "'126.17"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('D51').Value = "'1.197"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.05%  '
